# Logged Week 17 data and fixed Simulate_Season.py tiebreaking method
$wb = $excel.ActiveWorkbook

# --- OFF sheet: update Road ("R") row totals ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 463
$wsOff.Range("C3").Value = 306
$wsOff.Range("D3").Value = 120
$wsOff.Range("E3").Value = 56
$wsOff.Range("F3").Value = 8
$wsOff.Range("G3").Value = 5

# --- DEF sheet: update Road ("R") row totals ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 448
$wsDef.Range("C3").Value = 298
$wsDef.Range("D3").Value = 83
$wsDef.Range("E3").Value = 32
$wsDef.Range("F3").Value = 12
$wsDef.Range("G3").Value = 3
